$d = $word.ActiveDocument

$newText = "Os mapas de estrelas deste documento foron preparados por Jenik Hollan, CzechGlobe (http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)."

foreach ($p in $d.Paragraphs) {
    $pr = $p.Range
    if ($pr.Text -like "*Os mapas de*CzechGlobe*GaNight*") {
        # Paragraph.Range includes the trailing paragraph mark; exclude it so
        # we only touch the run content, not the paragraph break itself.
        $contentStart = $pr.Start
        $contentEnd = $pr.End - 1

        $contentRange = $d.Range($contentStart, $contentEnd)
        $contentRange.Text = ""

        $insertionPoint = $d.Range($contentStart, $contentStart)
        $insertionPoint.InsertAfter($newText)
        break
    }
}
